$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the return-rate assumption in B42 (20% -> 15%); all downstream
# pension-simulation formulas reference B$42 and recalculate automatically.
$ws.Range("B42").Value = 0.15

# Restore the view/selection state captured in the saved workbook.
$ws.Activate()
$ws.Range("O90").Select()
